$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Neo"
$ws.Range("B4").Value = "USA"
